$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.865.29"
$ws.Range("E2").Value = "  -2.31%  "
$ws.Range("D3").Value = "2.912.49"
$ws.Range("E3").Value = "  -3.27%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "2.910.80"
$ws.Range("E8").Value = "  -3.27%  "
$ws.Range("E9").Value = "  -3.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.68"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.144"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.43%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.446"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000224"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("D16").Value = "3.397.47"
$ws.Range("E16").Value = "  -3.21%  "
$ws.Range("D17").Value = "60.893.05"
$ws.Range("E17").Value = "  -2.19%  "
$ws.Range("E18").Value = "  -2.46%  "
$ws.Range("D19").Value = "2.912.58"
$ws.Range("E19").Value = "  -3.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "425.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.35%  "
$ws.Range("E22").Value = "  -2.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.92"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.78%  "
$ws.Range("E26").Value = "  -2.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.75%  "
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.02%  "
$ws.Range("E31").Value = "  -2.96%  "
$ws.Range("E32").Value = "  +3.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.62"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.90%  "
$ws.Range("E34").Value = "  -4.00%  "
$ws.Range("D35").Value = "0.0₃0845"
$ws.Range("E35").Value = "  -0.22%  "
$ws.Range("E36").Value = "  -1.46%  "
$ws.Range("E37").Value = "  -2.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.77"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.03"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.124"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.74"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.54"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.286"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.65%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "373.93"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0344"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.25%  "
$ws.Range("D47").Value = "2.648.80"
$ws.Range("E47").Value = "  -2.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.93"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.42%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.41"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.17%  "
